$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old literal/formula content in B29:E32 (now blank cells)
$ws.Range("B29:E32").ClearContents()

# Update the selection to match the new state
$ws.Range("B29:E32").Select()
